$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds crypto prices stored as literal text (values like "1.001" look
# numeric, so Excel would otherwise silently convert them to real numbers on
# assignment). Temporarily flip each touched D cell to a text format before
# writing the new price string, then drop back to the default/unstyled look so
# no visible formatting changes are introduced.
$dCells = @("D2", "D3", "D4", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D49", "D50", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "23.308.07"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "1.622.15"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "302.46"
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("D7").Value = "0.3751"
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("D8").Value = "51.51"
$ws.Range("E8").Value = "  -1.43%  "
$ws.Range("D9").Value = "0.3613"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").Value = "0.08128"
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").Value = "1.218"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "22.24"
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("D14").Value = "6.451"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").Value = "0.00001231"
$ws.Range("E15").Value = "  -2.94%  "
$ws.Range("D16").Value = "7.252"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").Value = "1.617.53"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "93.95"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").Value = "0.06925"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").Value = "17.47"
$ws.Range("E20").Value = "  -3.61%  "
$ws.Range("D21").Value = "6.521"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "12.46"
$ws.Range("E23").Value = "  -2.09%  "
$ws.Range("D24").Value = "23.301.66"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").Value = "2.467"
$ws.Range("E25").Value = "  +2.37%  "
$ws.Range("D26").Value = "3.073"
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("D27").Value = "21.10"
$ws.Range("E27").Value = "  -0.52%  "
$ws.Range("D28").Value = "150.55"
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").Value = "5.273"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").Value = "132.68"
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("D31").Value = "1.797.71"
$ws.Range("D32").Value = "6.711"
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").Value = "2.159"
$ws.Range("E33").Value = "  -5.50%  "
$ws.Range("D34").Value = "1.054"
$ws.Range("E34").Value = "  +10.49%  "
$ws.Range("D35").Value = "11.22"
$ws.Range("E35").Value = "  +8.61%  "
$ws.Range("D36").Value = "0.02741"
$ws.Range("E36").Value = "  -3.37%  "
$ws.Range("D37").Value = "0.08784"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "0.2473"
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("D39").Value = "0.07077"
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("D40").Value = "5.966"
$ws.Range("E40").Value = "  -1.43%  "
$ws.Range("D41").Value = "0.6952"
$ws.Range("D42").Value = "1.327"
$ws.Range("E42").Value = "  -3.53%  "
$ws.Range("D43").Value = "16.05"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D45").Value = "0.6439"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("E48").Value = "  -2.83%  "
$ws.Range("D49").Value = "0.07955"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "125.60"
$ws.Range("E50").Value = "  -2.12%  "
$ws.Range("D51").Value = "1.180"
$ws.Range("E51").Value = "  -1.38%  "

# Restore default (General/Normal) styling on the price cells we touched.
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
